$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = -5.889999999999998
$ws.Range("A9").Value = -20.55319999999999
$ws.Range("D11").Value = -8.562099999999999
$ws.Range("A18").Value = -22.95280000000001
$ws.Range("A20").Value = -22.08020000000002
